$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the IMEI column (old column E)
$ws.Columns("E").Delete()

# 2. Insert a new column for "mac" before the (now shifted) longitude column (I)
$ws.Columns("I").Insert()
